$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the last existing header cell (G1) to the new
# header cell (H1), matching the bold/border/center style used by the
# other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header label and the new data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
